$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Good Morning" cell to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8 to match the recorded selection state
$ws.Range("E8").Select()
